$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet -> Sheet1)
$ws.Name = "Sheet1"

# ---- Header row (row 1): new column headers ----
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Model Name"
$ws.Range("C1").Value = "Exact Precision (Micro Avg)"
$ws.Range("D1").Value = "Exact Recall (Micro Avg)"
$ws.Range("E1").Value = "Exact F1 Score (Micro Avg)"
$ws.Range("F1").Value = "Exact Precision (Macro Avg)"
$ws.Range("G1").Value = "Exact Recall (Macro Avg)"
$ws.Range("H1").Value = "Exact F1 Score (Macro Avg)"
$ws.Range("I1").Value = "Exact Precision (Weighted Avg)"
$ws.Range("J1").Value = "Exact Recall (Weighted Avg)"
$ws.Range("K1").Value = "Exact F1 Score (Weighted Avg)"
$ws.Range("L1").Value = "Partial Precision"
$ws.Range("M1").Value = "Partial Recall"
$ws.Range("N1").Value = "Partial F1 Score"
$ws.Range("O1").Value = "Partial TP"
$ws.Range("P1").Value = "Partial FP"
$ws.Range("Q1").Value = "Partial FN"
$ws.Range("R1").Value = "Support"
$ws.Range("S1").Value = "Accuracy"
$ws.Range("T1").Value = "Result Link"
$ws.Range("U1").Value = "Stats Link"
$ws.Range("V1").Value = "No of GPU Used"
$ws.Range("W1").Value = "Power Consumption"
$ws.Range("X1").Value = "Unnamed: 23"

# Header formatting: bold font, thin box border, centered/top aligned
$headerRange = $ws.Range("A1:X1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Row 2 data ----
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "09/11/2025"
$ws.Range("B2").Value = "Llama-3.1-8B-Instruct"
$ws.Range("C2").Value = 0.3815789473684211
$ws.Range("D2").Value = 0.2078853046594982
$ws.Range("E2").Value = 0.2691415313225058
$ws.Range("F2").Value = 0.2724752985391283
$ws.Range("G2").Value = 0.1143719833990947
$ws.Range("H2").Value = 0.153555716314337
$ws.Range("I2").Value = 0.4862831039335271
$ws.Range("J2").Value = 0.2078853046594982
$ws.Range("K2").Value = 0.2816885582180354
$ws.Range("L2").Value = 0.5133333333333333
$ws.Range("M2").Value = 0.276978417266187
$ws.Range("N2").Value = 0.3598130841121496
$ws.Range("O2").Value = 77
$ws.Range("P2").Value = 73
$ws.Range("Q2").Value = 201
$ws.Range("R2").Value = 279
$ws.Range("S2").Value = 0.9494178039313885
$ws.Range("T2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Llama-3.1-8B-Instruct_4_shot.txt"
$ws.Range("U2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Llama-3.1-8B-Instruct_4_shot.txt"
$ws.Range("V2").Value = "4 MLGPU"
$ws.Range("W2").Value = "0.042 kWh"
$ws.Range("X2").Value = 604

# ---- Row 3 data ----
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "09/12/2025"
$ws.Range("B3").Value = "Llama-3.1-8B-Instruct"
$ws.Range("C3").Value = 0.5029239766081871
$ws.Range("D3").Value = 0.2895622895622896
$ws.Range("E3").Value = 0.3675213675213675
$ws.Range("F3").Value = 0.5671276405298962
$ws.Range("G3").Value = 0.2904680106505578
$ws.Range("H3").Value = 0.3649174589104321
$ws.Range("I3").Value = 0.5722028737066331
$ws.Range("J3").Value = 0.2895622895622896
$ws.Range("K3").Value = 0.3696407221566667
$ws.Range("L3").Value = 0.5705882352941176
$ws.Range("M3").Value = 0.3277027027027027
$ws.Range("N3").Value = 0.4163090128755365
$ws.Range("O3").Value = 97
$ws.Range("P3").Value = 73
$ws.Range("Q3").Value = 199
$ws.Range("R3").Value = 297
$ws.Range("S3").Value = 0.9553025763930497
$ws.Range("T3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Llama-3.1-8B-Instruct_4_shot.txt"
$ws.Range("U3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Llama-3.1-8B-Instruct_4_shot.txt"
$ws.Range("V3").Value = "4 MLGPU"
$ws.Range("W3").Value = "0.025 kWh"
$ws.Range("X3").NumberFormat = "General"

Write-Output "done"
